$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new date row (A5) ---------------------------------------------
# The sheet currently has data through row 4 (A1:B4). The bug fix adds a new
# date value in A5 (serial 41264 -> 2012-12-21) formatted as a date so the
# SAX reader can correctly detect/parse it as a date cell.
$ws.Range("A5").Value = 41264
$ws.Range("A5").NumberFormat = "m/d/yyyy h:mm"

# --- Widen column A so the new date value is fully visible ------------------
$ws.Columns.Item(1).ColumnWidth = 18.3

# --- Move the active selection to the newly added cell ----------------------
$ws.Range("A5").Select() | Out-Null
